$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column D "Número Curso" with a sequential course number for the
# first few course rows (2-4). Existing columns A (cursosEncontrados),
# B (urlCursos) and C (Estatus) together with their data rows are unchanged.
$ws.Range("D1").Value = "Número Curso"
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(4, 4).Value = 4

# Widen column A and B so the long course names / urls are fully visible.
# (closest achievable values given the runtime's internal character-width
# quantization; targets are 65.42578125 and 29.5703125 "characters")
$ws.Columns.Item(1).ColumnWidth = 64.6
$ws.Columns.Item(2).ColumnWidth = 28.65

# Leave the selection on I5, matching the final state of the workbook.
$ws.Range("I5").Select()
